$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.426.17'
$ws.Range('E2').Value = '  -0.64%  '

$ws.Range('D3').Value = '1.725.09'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.55'
$ws.Range('E5').Value = '  -0.92%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.02%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4903'
$ws.Range('E7').Value = '  +1.61%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2606'
$ws.Range('E8').Value = '  -2.40%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06202'

$ws.Range('D10').Value = '1.720.51'
$ws.Range('E10').Value = '  -0.57%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06986'

$ws.Range('E12').Value = '  -0.73%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.539'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5995'
$ws.Range('E14').Value = '  -1.98%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.43'
$ws.Range('E15').Value = '  +0.33%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9998'
$ws.Range('E16').Value = '  -0.05%  '

$ws.Range('D17').Value = '26.420.94'
$ws.Range('E17').Value = '  -0.66%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9999'
$ws.Range('E18').Value = '  -0.06%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007228'
$ws.Range('E19').Value = '  +3.73%  '

$ws.Range('E20').Value = '  -1.57%  '

$ws.Range('D21').Value = '1.937.40'
$ws.Range('E21').Value = '  -0.96%  '

$ws.Range('E22').Value = '  -1.22%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.574'
$ws.Range('E23').Value = '  -2.51%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.150'
$ws.Range('E24').Value = '  -1.73%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.63'
$ws.Range('E25').Value = '  +0.26%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.28'
$ws.Range('E26').Value = '  -1.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.398'
$ws.Range('E27').Value = '  -0.76%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '106.91'
$ws.Range('E28').Value = '  -1.03%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.726'
$ws.Range('E29').Value = '  -2.25%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.952'
$ws.Range('E30').Value = '  -0.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08000'
$ws.Range('E31').Value = '  -0.12%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.681'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04502'
$ws.Range('E33').Value = '  -1.32%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9992'
$ws.Range('E34').Value = '  -0.03%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.602'
$ws.Range('E35').Value = '  -0.45%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.004'
$ws.Range('E36').Value = '  +0.14%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6266'
$ws.Range('E37').Value = '  -0.54%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9409'
$ws.Range('E38').Value = '  +4.81%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.387'
$ws.Range('E39').Value = '  +0.26%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.948'
$ws.Range('E40').Value = '  -3.47%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9998'
$ws.Range('E41').Value = '  -0.39%  '

$ws.Range('E42').Value = '  -1.09%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.62'
$ws.Range('E43').Value = '  -2.93%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.307'
$ws.Range('E44').Value = '  -2.15%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3856'
$ws.Range('E45').Value = '  -0.88%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.832'
$ws.Range('E46').Value = '  -4.04%  '

$ws.Range('E47').Value = '  -1.07%  '

$ws.Range('E48').Value = '  -0.50%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.748'
$ws.Range('E49').Value = '  -2.16%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '30.22'
$ws.Range('E50').Value = '  -1.21%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.237'
$ws.Range('E51').Value = '  -1.21%  '
